$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level value edits within rows 2-25 (no row shift) ---

# D2: -13.5 -> blank (missing)
$ws.Range("D2").Value = $null

# D5: blank -> -14.4
$ws.Range("D5").Value = -14.4

# C6: blank -> 15.1, D6: blank -> -14.2
$ws.Range("C6").Value = 15.1
$ws.Range("D6").Value = -14.2

# C8: 15.5 -> blank
$ws.Range("C8").Value = $null

# D9: -14.5 -> blank
$ws.Range("D9").Value = $null

# D10: -14.7 -> blank
$ws.Range("D10").Value = $null

# C12: blank -> 12.5
$ws.Range("C12").Value = 12.5

# C14: 14.4 -> blank
$ws.Range("C14").Value = $null

# C17: blank -> 11.2
$ws.Range("C17").Value = 11.2

# C18: blank -> 11.5
$ws.Range("C18").Value = 11.5

# C19: 13.2 -> blank
$ws.Range("C19").Value = $null

# C20: 12.5 -> blank
$ws.Range("C20").Value = $null

# C23: blank -> 12.2
$ws.Range("C23").Value = 12.2

# D24: blank -> -13.9
$ws.Range("D24").Value = -13.9

# --- Rows 26-33 get fully replaced data (rows 34-35 will be deleted) ---

$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = $null
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = $null
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = $null
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = $null
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = $null
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = $null
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

# --- Delete now-obsolete rows 34 and 35 ---
$ws.Range("A34:F35").Delete()
